$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 966 round (columns S-Z) line "라인" numbers, row 9
$ws.Range("T9").Value = 6
$ws.Range("U9").Value = 14
$ws.Range("V9").Value = 19
$ws.Range("W9").Value = 26
$ws.Range("X9").Value = 34
$ws.Range("Y9").Value = 45

# row 10
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 11
$ws.Range("V10").Value = 22
$ws.Range("W10").Value = 26
$ws.Range("X10").Value = 31
$ws.Range("Y10").Value = 45

# row 11
$ws.Range("T11").Value = 3
$ws.Range("U11").Value = 4
$ws.Range("V11").Value = 15
$ws.Range("W11").Value = 23
$ws.Range("X11").Value = 30
$ws.Range("Y11").Value = 45

# row 12
$ws.Range("T12").Value = 13
$ws.Range("U12").Value = 9
$ws.Range("V12").Value = 10
$ws.Range("W12").Value = 37
$ws.Range("X12").Value = 45
$ws.Range("Y12").Value = 39

# row 13
$ws.Range("T13").Value = 5
$ws.Range("U13").Value = 12
$ws.Range("V13").Value = 15
$ws.Range("W13").Value = 23
$ws.Range("X13").Value = 32
$ws.Range("Y13").Value = 43

# 968 round (columns J-Q) "결과" row 14
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 12
$ws.Range("N14").Value = 14
$ws.Range("O14").Value = 24
$ws.Range("P14").Value = 39
$ws.Range("Q14").Value = 33

# Update the selection to match the diff (active cell X19)
$ws.Range("X19").Select()
